# Update odds data in Sheet1 for Jogos_do_Dia_Betfair_Back_Lay_2025-10-07.xlsx
# Applies the per-cell value corrections captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Scunthorpe vs Morecambe)
$ws.Range("H4").Value = 7.4
$ws.Range("K4").Value = 5.6

# Row 5 (Kelty Hearts vs Alloa)
$ws.Range("F5").Value = 3.45
$ws.Range("G5").Value = 8.199999999999999
$ws.Range("H5").Value = 1.66
$ws.Range("I5").Value = 2.12
$ws.Range("J5").Value = 3.65
$ws.Range("L5").Value = 1.28
$ws.Range("P5").Value = 1.87
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 1.32
$ws.Range("S5").Value = 2.64
$ws.Range("V5").Value = 1.89

# Row 6 (Caernarfon Town vs The New Saints)
$ws.Range("F6").Value = 1.16

# Row 11 (Real Soacha Cundinamarca FC vs Orsomarso)
$ws.Range("F11").Value = 1.48
$ws.Range("H11").Value = 5.7
$ws.Range("J11").Value = 3.4
$ws.Range("K11").Value = 6.6
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 1.11
$ws.Range("O11").Value = 1.09
$ws.Range("P11").Value = 1.52
$ws.Range("Q11").Value = 2.16
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 1.03

# Row 12 (Amazonas FC vs Criciuma)
$ws.Range("F12").Value = 2.84
$ws.Range("G12").Value = 3.1
$ws.Range("H12").Value = 2.82
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 2.98
$ws.Range("K12").Value = 3.25
$ws.Range("N12").Value = 2.68
$ws.Range("O12").Value = 1.51
$ws.Range("P12").Value = 1.57
$ws.Range("Q12").Value = 2.44
$ws.Range("R12").Value = 1.21
$ws.Range("T12").Value = 2.04
$ws.Range("U12").Value = 1.84
$ws.Range("V12").Value = 1.47
$ws.Range("W12").Value = 1.47
$ws.Range("AI12").Value = 70
$ws.Range("AK12").Value = 980
$ws.Range("AL12").Value = 70
